# Add benchmark result row for "Atom N270" to the worksheet, right below
# the last existing row (row 45), following the same layout as the other
# rows in the "Tabulka1" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

$ws.Cells.Item($row, 1).Value = "Intel"        # A - Manufacturer
$ws.Cells.Item($row, 2).Value = 1              # B - #
$ws.Cells.Item($row, 3).Value = "Atom N270"    # C - CPU Model
$ws.Cells.Item($row, 4).Value = 2.5            # D - TDP (W)
$ws.Cells.Item($row, 5).Value = 1              # E - Cores
$ws.Cells.Item($row, 6).Value = 2              # F - SMT
$ws.Cells.Item($row, 7).Value = 1.6            # G - Base (GHz)
$ws.Cells.Item($row, 8).Value = 1.6            # H - Turbo (GHz)
$ws.Cells.Item($row, 9).Value = "x86-32"       # I - OS Architecture
# J (Power profile) intentionally left blank for this entry
$ws.Cells.Item($row, 11).Value = 1             # K - RAM Amount (GB)
$ws.Cells.Item($row, 12).Value = 1             # L - Channels
$ws.Cells.Item($row, 13).Value = "DDR2"        # M - Type
$ws.Cells.Item($row, 14).Value = 533           # N - Rating
$ws.Cells.Item($row, 15).Value = 40.65         # O -  26
$ws.Cells.Item($row, 16).Value = 83.55         # P -  27
$ws.Cells.Item($row, 17).Value = 221.53        # Q -  28
$ws.Cells.Item($row, 18).Value = "OOM"         # R -  29 (out of memory -> text, right aligned like other OOM cells)
$ws.Cells.Item($row, 18).HorizontalAlignment = -4152   # xlRight, matches the other "OOM" cells

# Move the active selection the way it ended up after the edit was made
$ws.Range("N47").Select() | Out-Null
